# 12/10/2017 CHITRA AND MAMATHA CHICK IN
#
# 1) Merge the two runs of the first (date) paragraph into a single run.
#    A "replace whole text with itself" via Find/Replace naturally
#    coalesces the two adjacent identically-formatted runs into one.
$d = $word.ActiveDocument

$d.Paragraphs.First.Range.Find.Execute(
    "Thu Sep 19 11:04:16 PDT 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "Thu Sep 19 11:04:16 PDT 2017", 2) | Out-Null

# 2) Append a brand-new purchase-detail record (TUE Oct 10 11:32:49 PDT
#    2017, DNR, CHOWCHOW) right after the existing "Amount balance"
#    line, before the trailing blank paragraphs.
#
# The existing trailing empty PlainText paragraph (non-bold) is used as
# the anchor: InsertParagraphBefore() on it always drops the new,
# still-empty paragraph immediately in front of it (and keeps the
# anchor's own formatting, i.e. plain / non-bold), so repeating the call
# builds the whole block in the correct forward order. NOTE: cached
# Paragraph objects and .Index do not renumber after mutation here, so
# every paragraph is re-fetched from Paragraphs.Item() using a freshly
# computed index instead of reusing references.
$countBefore = $d.Paragraphs.Count
$anchorIndex = $countBefore - 1
$anchor = $d.Paragraphs.Item($anchorIndex).Range

for ($i = 0; $i -lt 11; $i++) {
    $anchor.InsertParagraphBefore()
}

$base = $anchorIndex

# blank bold separator line
$d.Paragraphs.Item($base).Range.Font.Bold = $true

# date line
$d.Paragraphs.Item($base + 1).Range.InsertAfter("TUE Oct 10 11:32:49 PDT 2017")

# Person Name ... - DNR
$d.Paragraphs.Item($base + 2).Range.InsertAfter("Person Name`t`t`t`t- DNR")

# separator rule
$d.Paragraphs.Item($base + 3).Range.InsertAfter("---------------------------------------------------------------")

# Item Name ... - CHOWCHOW
$d.Paragraphs.Item($base + 4).Range.InsertAfter("Item Name`t`t`t`t- CHOWCHOW")

# Number of Pockets ... - 1
$d.Paragraphs.Item($base + 5).Range.InsertAfter("Number of Pockets`t`t`t- 1")

# Number of KGs ... - 76
$d.Paragraphs.Item($base + 6).Range.InsertAfter("Number of KGs`t`t`t- 76")

# Rate ... - 13
$d.Paragraphs.Item($base + 7).Range.InsertAfter("Rate`t`t`t`t`t- 13")

# Total Price ... - 988.0
$d.Paragraphs.Item($base + 8).Range.InsertAfter("Total Price`t`t`t`t- 988.0")

# Amount balance (bold) ... - 12024.0
$p20 = $d.Paragraphs.Item($base + 9)
$p20.Range.Font.Bold = $true
$p20.Range.InsertAfter("Amount balance`t`t`t- 12024.0")

# trailing blank (non-bold): $d.Paragraphs.Item($base + 10) left empty/plain on purpose

# trailing blank (bold), precedes the two pre-existing blank paragraphs
$d.Paragraphs.Item($base + 11).Range.Font.Bold = $true
